# Update the two-digit / one-digit division answers in the worksheet table.
# Each "old" value occurs exactly once in the document at the time it is
# searched for (verified against the full change-set), so a simple
# Find/Replace-all per pair, applied in this order, reproduces the diff
# without any cross-collisions between old/new values that coincide.

$d = $word.ActiveDocument

$replacements = @(
    @("85÷3=28, 1", "93÷3=31, 0"),
    @("25÷3=8, 1", "65÷3=21, 2"),
    @("97÷2=48, 1", "81÷8=10, 1"),
    @("15÷7=2, 1", "56÷4=14, 0"),
    @("16÷9=1, 7", "12÷3=4, 0"),
    @("77÷4=19, 1", "98÷9=10, 8"),
    @("20÷3=6, 2", "10÷8=1, 2"),
    @("76÷4=19, 0", "15÷3=5, 0"),
    @("76÷2=38, 0", "16÷9=1, 7"),
    @("96÷3=32, 0", "75÷4=18, 3"),
    @("44÷2=22, 0", "14÷2=7, 0"),
    @("35÷2=17, 1", "43÷3=14, 1"),
    @("82÷7=11, 5", "69÷4=17, 1"),
    @("27÷5=5, 2", "96÷5=19, 1"),
    @("34÷5=6, 4", "28÷8=3, 4"),
    @("59÷6=9, 5", "87÷2=43, 1"),
    @("23÷9=2, 5", "89÷9=9, 8"),
    @("16÷3=5, 1", "80÷5=16, 0"),
    @("34÷7=4, 6", "42÷9=4, 6"),
    @("82÷5=16, 2", "85÷9=9, 4"),
    @("79÷3=26, 1", "78÷2=39, 0"),
    @("16÷5=3, 1", "56÷8=7, 0"),
    @("57÷4=14, 1", "59÷6=9, 5"),
    @("65÷5=13, 0", "70÷3=23, 1"),
    @("56÷5=11, 1", "87÷6=14, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
